$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates (new "False positives" / "False negatives" values, and Notes) ---

# Row 4
$ws.Range("C4").Value = "na"
$ws.Range("D4").Value = "na"

# Row 5
$ws.Range("C5").Value = "na"
$ws.Range("D5").Value = "na"

# Row 6
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1

# Row 7
$ws.Range("C7").Value = "na"
$ws.Range("D7").Value = "na"

# Row 8
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 0
$ws.Range("F8").Value = "look to be actually 7 cells in the frame -- AP detected 2 false pos. too"

# Row 9
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = "Analyze particles based on bad mask"

# Row 10
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0

# Row 11
$ws.Range("C11").Value = "na"
$ws.Range("D11").Value = "na"

# --- Sheet view: freeze the header row and select E12 ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E12").Select() | Out-Null
